$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "服裝"
$ws.Range("B6").Value = 300
$ws.Range("C6").Value = "李湘菱"
$ws.Range("D6").Value = "other"

$ws.Range("A7:XFD7").Select()
